# "Changed Structure to dynamically loading File Reader Classes!"
#
# The "2nd Sheet" tab used to hold its own standalone placeholder rows
# (B/C columns = "1".."8", D/E always "Test"/blank). After switching the
# DataDriver example to load file-reader classes dynamically, that sheet's
# sample rows were replaced with the same data pattern already used on the
# "DataDriven" sheet (rows 2-9, columns A-E), and the leftover placeholder
# strings became unused and were dropped from the shared-string table.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item('2nd Sheet')

# --- Row 2 --------------------------------------------------------------
$ws2.Cells.Item(2,1).Value = 'Right user empty pass'
$ws2.Cells.Item(2,2).Value = 'demo'
$ws2.Cells.Item(2,3).Value = '${EMPTY}'
# D2 must stay a genuine number (1), not text, even though the column's
# style carries a Text ("@") number format - flip to General just long
# enough to store the number, then restore the original Text format.
$d2 = $ws2.Cells.Item(2,4)
$d2.NumberFormat = 'General'
$d2.Value = 1
$d2.NumberFormat = '@'
$ws2.Cells.Item(2,5).Value = 'This is a test case documentation of the first one.'

# --- Row 3 --------------------------------------------------------------
$ws2.Cells.Item(3,1).Value = 'Right user wrong pass'
$ws2.Cells.Item(3,2).Value = 'demo'
$ws2.Cells.Item(3,3).Value = 'FooBar'
$ws2.Cells.Item(3,4).Value = '2,3,foo'
$ws2.Cells.Item(3,5).Value = 'This test case has the Tags 2,3 and foo'

# --- Row 4 --------------------------------------------------------------
$ws2.Cells.Item(4,2).Value = '${EMPTY}'
$ws2.Cells.Item(4,3).Value = 'mode'
$ws2.Cells.Item(4,4).Value = '1,2,3,4'
$ws2.Cells.Item(4,5).Value = 'This test case has a generated name based on template name.'

# --- Row 5 --------------------------------------------------------------
$ws2.Cells.Item(5,2).Value = '${EMPTY}'
$ws2.Cells.Item(5,3).Value = '${EMPTY}'

# --- Row 6 --------------------------------------------------------------
$ws2.Cells.Item(6,2).Value = '${EMPTY}'
$ws2.Cells.Item(6,3).Value = 'FooBar'
$ws2.Cells.Item(6,4).Value = 'foo'

# --- Row 7 --------------------------------------------------------------
$ws2.Cells.Item(7,2).Value = 'FooBar'
$ws2.Cells.Item(7,3).Value = 'mode'
$ws2.Cells.Item(7,4).Value = 'foo'

# --- Row 8 --------------------------------------------------------------
$ws2.Cells.Item(8,2).Value = 'FooBar'
$ws2.Cells.Item(8,3).Value = '${EMPTY}'
$ws2.Cells.Item(8,4).Value = 'foo'

# --- Row 9 --------------------------------------------------------------
$ws2.Cells.Item(9,2).Value = 'FooBar'
$ws2.Cells.Item(9,3).Value = 'FooBar'
$ws2.Cells.Item(9,4).Value = 'foo'

# The author last had G7 selected on "2nd Sheet" when the file was saved.
$ws2.Activate()
$ws2.Range('G7').Select()
